# "more datacamp and disk_savvy data"
#
# A new DataCamp course entry ("Software Engineering Principles in
# Python") is being logged. It is inserted at row 35 (where the course
# list was being edited), and the row that used to live at 35
# ("Software Engineering for Data Scientists in Python", with its
# A:E rating columns) is pushed down to the newly-created row 40.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the course that used to occupy row 35 down to row 40,
#     preserving its rating values and the highlighted-title font
#     color used for course-name cells in this block ---
$ws.Range("A40").Value = $ws.Range("A35").Value()
$ws.Range("B40").Value = $ws.Range("B35").Value()
$ws.Range("C40").Value = $ws.Range("C35").Value()
$ws.Range("D40").Value = $ws.Range("D35").Value()
$ws.Range("E40").Value = $ws.Range("E35").Value()
$ws.Range("A40").Font.Color = $ws.Range("A35").Font.Color()

# --- Clear the old row 35 rating values (B:E) - the new entry only
#     has a single score, in column F ---
$ws.Range("B35:E35").Clear()

# --- Write the new course entry into row 35 ---
$ws.Range("A35").Value = "Software Engineering Principles in Python"
$ws.Range("F35").Value = 2

# --- Leave the selection where the user's cursor ended up next ---
$ws.Range("A36").Select()
